$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns D/E content added alongside existing rows 25-28 (Otomatik Sarj sub-table)
$ws.Range("D25").Value = "Otomatik Şarj"
$ws.Range("E25").Value = "Float ve boost testleri tamamlanmış varsayılıyor."
$ws.Range("E26").Value = "Manuel olarak float şarj moduna geç."
$ws.Range("E27").Value = "Voltaj ve akımı ön panelden oku, ölçüm cihazı ile ölç ve yaz."
$ws.Range("E28").Value = "Otomatik şarja al."

# Row 29 heading changes from "Otomatik Şarj" to "Zamanlı Şarj" (new Timed Charge section begins)
$ws.Range("A29").Value = "Zamanlı Şarj"

# New rows 31-33 for the Zamanlı Şarj (Timed Charge) procedure steps
$ws.Range("B31").Value = "Zamanlı şarj moduna geç."
$ws.Range("B32").Value = "Zamanı 1dk ya ayarla."
$ws.Range("B33").Value = "Ana menüye gel."

# Old row 33 content ("Oto şarj test") in column A is no longer present in the final layout
$ws.Range("A33").ClearContents()

$ws.Range("B36").Select()
